$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells I1 (I0) and J1 (IF), styled like the other header cells ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- New data columns I (I0) and J (IF) for rows 2-54 ---
$data = @(
    @(6,7),
    @(9,9),
    @(7,7),
    @(7,7),
    @(7,9),
    @(7,7),
    @(8,8),
    @(6,9),
    @(8,10),
    @(5,6),
    @(9,9),
    @(8,8),
    @(9,9),
    @(7,9),
    @(9,9),
    @(7,8),
    @(9,9),
    @(6,8),
    @(7,8),
    @(8,8),
    @(1,1),
    @(8,8),
    @(9,9),
    @(7,8),
    @(8,8),
    @(9,9),
    @(1,1),
    @(9,9),
    @(5,6),
    @(1,1),
    @(8,8),
    @(6,6),
    @(7,7),
    @(9,9),
    @(6,7),
    @(7,8),
    @(8,9),
    @(8,8),
    @(1,4),
    @(7,8),
    @(8,9),
    @(8,9),
    @(6,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(8,9),
    @(4,4),
    @(5,6),
    @(3,3),
    @(5,6),
    @(7,7),
    @(6,6)
)

for ($k = 0; $k -lt $data.Count; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $data[$k][0]
    $ws.Cells.Item($row, 10).Value = $data[$k][1]
}
